$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 15
$ws.Range("B1").Value = 3.063336372375488
$ws.Range("C1").Value = 2.713537931442261
$ws.Range("D1").Value = 2.959021329879761
$ws.Range("E1").Value = 15
